# Generate Report for Archive
# The localization status has moved on from "Ready for handoff" to
# "In Translation" for this file, across every place that status is
# reported: the Overview roll-up sheet (per-language columns) and each
# per-language detail sheet's "Status" column. Re-running the report
# generator also re-flows (auto-sizes) the now-narrower Status columns.

$wb = $excel.ActiveWorkbook

# --- Overview sheet -------------------------------------------------
# Columns E (zh-cn) and F (de-de) hold the current status per language.
$ws1 = $wb.Worksheets.Item("Overview")
$ws1.Range("E2").Value = "In Translation"
$ws1.Range("F2").Value = "In Translation"
$ws1.Columns.Item(5).ColumnWidth = 12.45
$ws1.Columns.Item(6).ColumnWidth = 12.45

# --- zh-cn detail sheet ----------------------------------------------
# Column C is the "Status" column.
$ws2 = $wb.Worksheets.Item("zh-cn")
$ws2.Range("C2").Value = "In Translation"
$ws2.Columns.Item(3).ColumnWidth = 12.45

# --- de-de detail sheet ----------------------------------------------
# Column C is the "Status" column.
$ws3 = $wb.Worksheets.Item("de-de")
$ws3.Range("C2").Value = "In Translation"
$ws3.Columns.Item(3).ColumnWidth = 12.45
